# Brynhildr_Profits workbook refresh (scheduled-runner price update)
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to the latest market data.
# A few rows also gain or lose a profit cell (NQ/HQ wasn't craftable before/now is).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 19813
$ws.Range("J106").Value = 19499.5
$ws.Range("L106").Value = 19499.5
$ws.Range("N106").Value = -20761.5
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 3401.647
$ws.Range("J112").Value = 3347.3572
$ws.Range("L112").Value = 10042.0716
$ws.Range("N112").Value = -12258.0716
# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 4238.1113
$ws.Range("J113").Value = 5036
$ws.Range("L113").Value = 5036
$ws.Range("N113").Value = -11544
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 5700.558
$ws.Range("I132").Value = 6155.5386
$ws.Range("K132").Value = 18466.6158
$ws.Range("M132").Value = -15936.6158
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 5510.857
$ws.Range("I135").Value = 264.5
$ws.Range("K135").Value = 2380.5
$ws.Range("M135").Value = 154.5
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 20005328
$ws.Range("I137").Value = 29413102
$ws.Range("J137").Value = 13809.875
$ws.Range("K137").Value = 88239306
$ws.Range("L137").Value = 41429.625
$ws.Range("M137").Value = -88236756
$ws.Range("N137").Value = -46529.625
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1775.3549
$ws.Range("I2").Value = 1571.52
$ws.Range("J2").Value = 2624.6667
$ws.Range("K2").Value = 1571.52
$ws.Range("L2").Value = 2624.6667
$ws.Range("M2").Value = -1458.52
$ws.Range("N2").Value = -2850.6667
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1775.3549
$ws.Range("I116").Value = 1571.52
$ws.Range("J116").Value = 2624.6667
$ws.Range("K116").Value = 1571.52
$ws.Range("L116").Value = 2624.6667
$ws.Range("M116").Value = 722.48
$ws.Range("N116").Value = -7212.6667
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 6424.156
$ws.Range("I132").Value = 4292.154
$ws.Range("K132").Value = 12876.462
$ws.Range("M132").Value = -10346.462
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1775.3549
$ws.Range("I3").Value = 1571.52
$ws.Range("J3").Value = 2624.6667
$ws.Range("K3").Value = 1571.52
$ws.Range("L3").Value = 2624.6667
$ws.Range("M3").Value = -1457.52
$ws.Range("N3").Value = -2852.6667
# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 433.33334
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -227
$ws.Range("N22").Value = -846
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1779.5
$ws.Range("I105").Value = 1629.138
$ws.Range("J105").Value = 2175.9092
$ws.Range("K105").Value = 1629.138
$ws.Range("L105").Value = 2175.9092
$ws.Range("M105").Value = 117.8620000000001
$ws.Range("N105").Value = -5669.9092
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 6175562
$ws.Range("I134").Value = 2899.625
$ws.Range("J134").Value = 55556860
$ws.Range("K134").Value = 8698.875
$ws.Range("L134").Value = 166670580
$ws.Range("M134").Value = -6163.875
$ws.Range("N134").Value = -166675650
$ws = $wb.Worksheets.Item("CRP")
# Row 18: Life's a Stitch / Ash Spinning Wheel
$ws.Range("H18").Value = 41497.25
$ws.Range("I18").Value = 20989
$ws.Range("J18").Value = 48333.332
$ws.Range("K18").Value = 20989
$ws.Range("L18").Value = 48333.332
$ws.Range("M18").Value = -20759
$ws.Range("N18").Value = -48793.332
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3078.9033
$ws.Range("I132").Value = 3106.0833
$ws.Range("K132").Value = 9318.249899999999
$ws.Range("M132").Value = -6788.249899999999
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1659.0938
$ws.Range("I134").Value = 1606.2858
$ws.Range("J134").Value = 2028.75
$ws.Range("K134").Value = 4818.857400000001
$ws.Range("L134").Value = 6086.25
$ws.Range("M134").Value = -2283.857400000001
$ws.Range("N134").Value = -11156.25
$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 10714388
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents() | Out-Null
# Row 18: Fisher of Men / Salt Cod
$ws.Range("H18").Value = 1767.7142
$ws.Range("I18").Value = 1102
$ws.Range("J18").Value = 2655.3333
$ws.Range("K18").Value = 3306
$ws.Range("L18").Value = 7965.999899999999
$ws.Range("M18").Value = -3137
$ws.Range("N18").Value = -8303.999899999999
# Row 86: Let's Not Get Sappy / Birch Syrup
$ws.Range("H86").Value = 821.2857
$ws.Range("I86").Value = 850
$ws.Range("J86").Value = 809.8
$ws.Range("K86").Value = 2550
$ws.Range("L86").Value = 2429.4
$ws.Range("M86").Value = -1364
$ws.Range("N86").Value = -4801.4
# Row 89: Luxury Spillover (L) / Birch Syrup
$ws.Range("H89").Value = 821.2857
$ws.Range("I89").Value = 850
$ws.Range("J89").Value = 809.8
$ws.Range("K89").Value = 7650
$ws.Range("L89").Value = 7288.2
$ws.Range("M89").Value = -1722
$ws.Range("N89").Value = -19144.2
# Row 114: One Last Meal / Mushroom Saute
$ws.Range("H114").Value = 5957.148
$ws.Range("J114").Value = 7843.4
$ws.Range("L114").Value = 23530.2
$ws.Range("N114").Value = -30038.2
# Row 123: Topping Up the Pot / Zurek
$ws.Range("H123").Value = 14916.667
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents() | Out-Null
# Row 124: Bobbing for Compliments / Island Miq'abob
$ws.Range("H124").Value = 27885.285
$ws.Range("I124").Value = 60098.5
$ws.Range("K124").Value = 180295.5
$ws.Range("M124").Value = -175385.5
# Row 125: At Any Temperature / Borscht
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents() | Out-Null
# Row 126: Imperial Palate / Glory Be Soup
$ws.Range("H126").Value = 12613.125
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 12613.125
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents() | Out-Null
$ws.Range("M126").Value = 37839.375
$ws.Range("N126").Value = -47719.375
# Row 130: Blast from the Pasta / The Noodles of Elpis
$ws.Range("H130").Value = 11032.875
$ws.Range("I130").Value = 4015
$ws.Range("K130").Value = 12045
$ws.Range("M130").Value = -7025
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 4788.4
$ws.Range("J131").Value = 5164.9443
$ws.Range("L131").Value = 15494.8329
$ws.Range("N131").Value = -25574.8329
# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 7417.05
$ws.Range("I133").Value = 3810.0715
$ws.Range("K133").Value = 11430.2145
$ws.Range("M133").Value = -6370.2145
# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 8107.75
$ws.Range("I134").Value = 1137.3334
$ws.Range("J134").Value = 12290
$ws.Range("K134").Value = 3412.0002
$ws.Range("L134").Value = 36870
$ws.Range("M134").Value = 1657.9998
$ws.Range("N134").Value = -47010
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 7790.1875
$ws.Range("I137").Value = 2405.5
$ws.Range("J137").Value = 13174.875
$ws.Range("K137").Value = 7216.5
$ws.Range("L137").Value = 39524.625
$ws.Range("M137").Value = -2116.5
$ws.Range("N137").Value = -49724.625
# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 25809.12
$ws.Range("I138").Value = 50461.145
$ws.Range("J138").Value = 16222.223
$ws.Range("K138").Value = 151383.435
$ws.Range("L138").Value = 48666.669
$ws.Range("M138").Value = -146243.435
$ws.Range("N138").Value = -58946.669
# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 8432.066000000001
$ws.Range("J139").Value = 11766.667
$ws.Range("L139").Value = 35300.001
$ws.Range("N139").Value = -45580.001
# Row 141: Ocean Explosion / Acqua Pazza
$ws.Range("H141").Value = 14381.444
$ws.Range("I141").Value = 9800
$ws.Range("J141").Value = 16672.166
$ws.Range("K141").Value = 29400
$ws.Range("L141").Value = 50016.49800000001
$ws.Range("M141").Value = -24220
$ws.Range("N141").Value = -60376.49800000001
$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3373.1738
$ws.Range("I122").Value = 3702.5
$ws.Range("J122").Value = 2860.889
$ws.Range("K122").Value = 11107.5
$ws.Range("L122").Value = 8582.667000000001
$ws.Range("M122").Value = -8657.5
$ws.Range("N122").Value = -13482.667
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 2597.8262
$ws.Range("I126").Value = 2409.6365
$ws.Range("K126").Value = 7228.9095
$ws.Range("M126").Value = -4758.9095
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 18928.143
$ws.Range("I132").Value = 22500.6
$ws.Range("K132").Value = 67501.79999999999
$ws.Range("M132").Value = -64971.79999999999
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 3605.6191
$ws.Range("I22").Value = 2433
$ws.Range("J22").Value = 3881.5293
$ws.Range("K22").Value = 2433
$ws.Range("L22").Value = 3881.5293
$ws.Range("M22").Value = -2138
$ws.Range("N22").Value = -4471.5293
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 3605.6191
$ws.Range("I27").Value = 2433
$ws.Range("J27").Value = 3881.5293
$ws.Range("K27").Value = 2433
$ws.Range("L27").Value = 3881.5293
$ws.Range("M27").Value = -2326
$ws.Range("N27").Value = -4095.5293
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 4810.3125
$ws.Range("I40").Value = 4433.615
$ws.Range("K40").Value = 4433.615
$ws.Range("M40").Value = -4297.615
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 4403.273
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 4893.263
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 4893.263
$ws.Range("M46").Value = -1112
$ws.Range("N46").Value = -5269.263
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1623.6154
$ws.Range("I93").Value = 928.1667
$ws.Range("K93").Value = 928.1667
$ws.Range("M93").Value = 319.8333
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 8619.75
$ws.Range("I122").Value = 4504
$ws.Range("K122").Value = 13512
$ws.Range("M122").Value = -11062
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3925279
$ws.Range("I132").Value = 4765203
$ws.Range("K132").Value = 14295609
$ws.Range("M132").Value = -14293079
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 5436777
$ws.Range("I136").Value = 3126842.5
$ws.Range("J136").Value = 20836338
$ws.Range("K136").Value = 9380527.5
$ws.Range("L136").Value = 62509014
$ws.Range("M136").Value = -9377977.5
$ws.Range("N136").Value = -62514114
$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 2458.3635
$ws.Range("I107").Value = 1169.35
$ws.Range("J107").Value = 4441.4614
$ws.Range("K107").Value = 3508.05
$ws.Range("L107").Value = 13324.3842
$ws.Range("M107").Value = -1588.05
$ws.Range("N107").Value = -17164.3842
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 375.2
$ws.Range("I113").Value = 180
$ws.Range("J113").Value = 424
$ws.Range("K113").Value = 540
$ws.Range("L113").Value = 1272
$ws.Range("M113").Value = 1630
$ws.Range("N113").Value = -5612
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 38158.562
$ws.Range("I122").Value = 2713.577
$ws.Range("K122").Value = 8140.731000000001
$ws.Range("M122").Value = -5690.731000000001
